# Applies the "all tasks completed except for command shell" edit to the
# session-based exploratory testing sheet.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) "Session length: long (3-4 hours)" -> "Session length: short (1-3 hours)"
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "Session length: long (3-4 hours)", $false, $false, $false, $false,
    $false, $true, 1, $false,
    "Session length: short (1-3 hours)", 2) | Out-Null

# ---------------------------------------------------------------------
# 2) "MacBook Pro 15” Majave v 10.14" -> same text, but "Majave" becomes
#    its own run (three runs total, no visible text change).
# ---------------------------------------------------------------------
$macRng = $d.Content
$macRng.Find.Execute(
    "Majave", $false, $false, $false, $false, $false, $true, 1, $false,
    "", 0) | Out-Null
# Toggling a character property and reverting it forces Word to split the
# found text into its own run without altering the visible formatting.
$macRng.Bold = 1
$macRng.Bold = 0

# ---------------------------------------------------------------------
# 3) "Views:" -> "Views / Flows:"
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "Views:", $false, $false, $false, $false, $false, $true, 1, $false,
    "Views / Flows:", 2) | Out-Null

# ---------------------------------------------------------------------
# 4) Insert a new bullet "User can log in successfully after waiting
#    appropriate time (~12 hrs)" right before "Password change flow
#    successful" (same list / indent formatting as its neighbours).
# ---------------------------------------------------------------------
$insRng = $d.Content
$insRng.Find.Execute(
    "Password change flow successful", $false, $false, $false, $false,
    $false, $true, 1, $false, "", 0) | Out-Null
$targetPara = $insRng.Paragraphs(1)
$targetPara.Range.InsertParagraphBefore()
# After the insert, $targetPara now refers to the freshly created (empty)
# paragraph that sits before the original "Password change flow
# successful" paragraph; it already inherited the surrounding bullet
# formatting, so we just need to give it its text.
$targetPara.Range.Text = "User can log in successfully after waiting appropriate time (~12 hrs)"

# ---------------------------------------------------------------------
# 5) Split "...results in 500 server error. " so "error" becomes its own
#    run (three runs total, no visible text change). Scope the Find to
#    the bug paragraph so we don't match the unrelated "error" earlier
#    in the document (re-materialise the paragraph bounds via
#    Document.Range so the scoped Find actually respects them).
# ---------------------------------------------------------------------
$bugRng = $d.Content
$bugRng.Find.Execute(
    "Bug# X: Selecting", $false, $false, $false, $false, $false, $true,
    1, $false, "", 0) | Out-Null
$bugPara = $bugRng.Paragraphs(1)
$errRng = $d.Range($bugPara.Range.Start, $bugPara.Range.End)
$errRng.Find.Execute(
    "error", $false, $false, $false, $false, $false, $true, 1, $false,
    "", 0) | Out-Null
$errRng.Bold = 1
$errRng.Bold = 0

Write-Output "edits applied"
